# B6-PowerPoint.pptx edit:
#  1) Re-style the three tables (slides 14, 15, 16) from the default
#     "Table_0" style to the built-in "Medium Style 2 - Accent 1" style.
#  2) Swap the deck's colour theme ("Integral" / Red-Violet) for the
#     stock "Office Theme" / Office colours (the Notes Master keeps the
#     Integral colours, matching the pre-edit Slide Master).

$p = $ppt.ActivePresentation

# -- 1. Table styles -------------------------------------------------
$newTableStyleId = "{54617B98-30B8-496D-86AC-CA3407425174}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# -- 2. Theme colours --------------------------------------------------
function ToBgr([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme (the scheme that now becomes active on the
# Slide Master, in ThemeColorScheme.Colors(1..12) order).
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = ToBgr($officeColors[$i - 1])
}
